$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new date column AV (column 48) with header "14-ago", matching the
# style used by the rest of the header row (AU1, etc.).
$ws.Cells.Item(1, 48).Value = "14-ago"
$ws.Cells.Item(1, 48).NumberFormat = "@"

# Fill in the new values for AV2:AV11, matching the number format/alignment
# used by the rest of the data column (AU2:AU11).
$ws.Range("AV2:AV11").NumberFormat = "0"
$ws.Range("AV2:AV11").HorizontalAlignment = -4108

$ws.Cells.Item(2, 48).Value = 16
$ws.Cells.Item(3, 48).Value = 15
$ws.Cells.Item(4, 48).Value = 11
$ws.Cells.Item(5, 48).Value = 11
$ws.Cells.Item(6, 48).Value = 16
$ws.Cells.Item(7, 48).Value = 11
$ws.Cells.Item(8, 48).Value = 12
$ws.Cells.Item(9, 48).Value = 21
$ws.Cells.Item(10, 48).Value = 16
$ws.Cells.Item(11, 48).Value = 29

# Update selection to mirror the recorded last-selected cell
$ws.Range("AY9").Select()
